# UndoRedoSequenceDiagram.pptx - rename the "AddressBook" sample domain to
# "FinanceTracker" across the sequence-diagram shapes on slide 1, and
# resize/reposition the affected label textboxes to match their new
# (longer) text.

$EMU_PER_PT = 12700.0

# Shape.Left/Top/Width/Height round-trip through a lower-precision
# (single/float) representation, so a naive "emu / 12700" assignment can
# land 1-2 EMU off target. Converge on the exact EMU value by re-reading
# the property and correcting for the residual error.
function Set-ExactLeft($shape, $targetEmu) {
    $pt = $targetEmu / $EMU_PER_PT
    $bestDiff = 999999999
    $bestPt = $pt
    for ($i = 0; $i -lt 40; $i++) {
        $shape.Left = $pt
        $gotExact = $shape.Left * $EMU_PER_PT
        $got = [math]::Round($gotExact)
        $diff = [math]::Abs($got - $targetEmu)
        if ($diff -lt $bestDiff) {
            $bestDiff = $diff
            $bestPt = $pt
        }
        if ($got -eq $targetEmu) { break }
        $err = $targetEmu - $gotExact
        $pt = $pt + ($err / $EMU_PER_PT)
    }
    $shape.Left = $bestPt
}
function Set-ExactTop($shape, $targetEmu) {
    $pt = $targetEmu / $EMU_PER_PT
    $bestDiff = 999999999
    $bestPt = $pt
    for ($i = 0; $i -lt 40; $i++) {
        $shape.Top = $pt
        $gotExact = $shape.Top * $EMU_PER_PT
        $got = [math]::Round($gotExact)
        $diff = [math]::Abs($got - $targetEmu)
        if ($diff -lt $bestDiff) {
            $bestDiff = $diff
            $bestPt = $pt
        }
        if ($got -eq $targetEmu) { break }
        $err = $targetEmu - $gotExact
        $pt = $pt + ($err / $EMU_PER_PT)
    }
    $shape.Top = $bestPt
}
function Set-ExactWidth($shape, $targetEmu) {
    $pt = $targetEmu / $EMU_PER_PT
    $bestDiff = 999999999
    $bestPt = $pt
    for ($i = 0; $i -lt 40; $i++) {
        $shape.Width = $pt
        $gotExact = $shape.Width * $EMU_PER_PT
        $got = [math]::Round($gotExact)
        $diff = [math]::Abs($got - $targetEmu)
        if ($diff -lt $bestDiff) {
            $bestDiff = $diff
            $bestPt = $pt
        }
        if ($got -eq $targetEmu) { break }
        $err = $targetEmu - $gotExact
        $pt = $pt + ($err / $EMU_PER_PT)
    }
    $shape.Width = $bestPt
}
function Set-ExactHeight($shape, $targetEmu) {
    $pt = $targetEmu / $EMU_PER_PT
    $bestDiff = 999999999
    $bestPt = $pt
    for ($i = 0; $i -lt 40; $i++) {
        $shape.Height = $pt
        $gotExact = $shape.Height * $EMU_PER_PT
        $got = [math]::Round($gotExact)
        $diff = [math]::Abs($got - $targetEmu)
        if ($diff -lt $bestDiff) {
            $bestDiff = $diff
            $bestPt = $pt
        }
        if ($got -eq $targetEmu) { break }
        $err = $targetEmu - $gotExact
        $pt = $pt + ($err / $EMU_PER_PT)
    }
    $shape.Height = $bestPt
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Lifeline header ":AddressBookParser" (two paragraphs: ":Address" /
#    "BookParser") becomes a single paragraph ":FinanceTrackerParser"
#    split into two runs (":" and "FinanceTrackerParser").
# ---------------------------------------------------------------------
$shpParser = $s.Shapes.Item(6)
$trParser = $shpParser.TextFrame.TextRange
$trParser.Text = ":FinanceTrackerParser"
$runParser = $trParser.Characters(2, "FinanceTrackerParser".Length)
$runParser.Text = "FinanceTrackerParser"

# ---------------------------------------------------------------------
# 2) "undoAddressBook()" label -> "undoFinanceTracker()" + reposition/resize.
# ---------------------------------------------------------------------
$shpUndo = $s.Shapes.Item(19)
$trUndo = $shpUndo.TextFrame.TextRange
$fullUndo = $trUndo.Text
$idxUndo = $fullUndo.IndexOf("AddressBook")
$cUndo = $trUndo.Characters($idxUndo + 1, "AddressBook".Length)
$cUndo.Text = "FinanceTracker"

Set-ExactLeft $shpUndo 5562600
Set-ExactTop $shpUndo 2748246
Set-ExactWidth $shpUndo 1502029
Set-ExactHeight $shpUndo 184666

# ---------------------------------------------------------------------
# 3) ":VersionedAddressBook" -> ":VersionedFinanceTracker" + resize (width).
# ---------------------------------------------------------------------
$shpVersioned = $s.Shapes.Item(23)
$trVersioned = $shpVersioned.TextFrame.TextRange
$fullVersioned = $trVersioned.Text
$idxVersioned = $fullVersioned.IndexOf("VersionedAddressBook")
$cVersioned = $trVersioned.Characters($idxVersioned + 1, "VersionedAddressBook".Length)
$cVersioned.Text = "VersionedFinanceTracker"

Set-ExactLeft $shpVersioned 7497155
Set-ExactTop $shpVersioned 2568606
Set-ExactWidth $shpVersioned 2328686
Set-ExactHeight $shpVersioned 335427

# ---------------------------------------------------------------------
# 4) "resetData(ReadOnlyAddressBook)" -> "resetData(ReadOnlyFinanceTracker)"
#    + reposition/resize.
# ---------------------------------------------------------------------
$shpReadOnly = $s.Shapes.Item(35)
$trReadOnly = $shpReadOnly.TextFrame.TextRange
$fullReadOnly = $trReadOnly.Text
$idxReadOnly = $fullReadOnly.IndexOf("ReadOnlyAddressBook")
$cReadOnly = $trReadOnly.Characters($idxReadOnly + 1, "ReadOnlyAddressBook".Length)
$cReadOnly.Text = "ReadOnlyFinanceTracker"

Set-ExactLeft $shpReadOnly 8534400
Set-ExactTop $shpReadOnly 3320534
Set-ExactWidth $shpReadOnly 2328686
Set-ExactHeight $shpReadOnly 184666
